$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = 40
$ws1.Range("F3").Value = 8519
$ws1.Range("F4").Value = 6227
$ws1.Range("F5").Value = 545
$ws1.Range("F6").Value = 117
$ws1.Range("F9").Value = 332
$ws1.Range("F10").Value = 1214

# Sheet "全部类型" (index 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = 40
$ws4.Range("F3").Value = 8519
$ws4.Range("F4").Value = 6227
$ws4.Range("F5").Value = 545
$ws4.Range("F6").Value = 117
$ws4.Range("F9").Value = 332
$ws4.Range("F14").Value = 1214
